$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.298.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.63%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.094.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.37%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.663"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.50%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.51"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +22.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.79"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.374"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0743"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.34%  "

$ws.Range("E12").Value = "  +8.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.400.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.837"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.097.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.250.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +14.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0847"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "241.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.97%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("E25").Value = "  +2.11%  "

$ws.Range("E26").Value = "  +4.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.53%  "

$ws.Range("E30").Value = "  +2.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +27.48%  "

$ws.Range("E33").Value = "  +4.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0617"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0911"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.68%  "

$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.87%  "

$ws.Range("E40").Value = "  +1.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +14.91%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0226"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.89%  "

$ws.Range("E43").Value = "  +5.92%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.81%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0930"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +106.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.324.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.47%  "

$ws.Range("E50").Value = "  +14.94%  "

$ws.Range("E51").Value = "  +6.50%  "
